$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RoleListingPageAttributes")

# New column G header + icon path values (renaming/updating shared strings in place)
$ws.Range("G1").Value = 'icon_value'
$ws.Range("G2").Value = 'path("M 11.13 8.524 C 11.47 8.254 11.495 7.896 11.27 7.419 A 19.186 19.186 0 0 1 10.46 5.445 A 3.194 3.194 0 0 1 10.414 3.553 C 10.663 2.631 11.37 2.013 12.348 1.883 C 13.057 1.789 13.733 1.859 14.335 2.28 C 15.233 2.908 15.4 3.811 15.295 4.79 C 15.199 5.698 14.858 6.54 14.419 7.344 C 14.044 8.032 14.185 8.469 14.931 8.753 C 15.418 8.938 15.935 9.047 16.438 9.191 C 16.645 9.251 16.859 9.294 17.055 9.377 C 17.435 9.539 17.716 9.797 17.763 10.225 C 17.813 10.675 17.954 11.53 18 12 H 12.72 C 12.7 11.577 12.681 11.265 12.662 11.064 C 12.623 10.651 12.598 10.264 12.569 9.999 C 12.484 9.22 12.052 8.843 11.229 8.562 A 1.79 1.79 0 0 1 11.13 8.524 Z M 0 12 C 0.09 10.737 0.097 10.387 0.298 9.767 C 0.423 9.381 0.743 9.131 1.135 9.012 C 1.693 8.843 2.265 8.719 2.825 8.557 A 6.36 6.36 0 0 0 3.805 8.207 C 4.508 7.87 4.625 7.382 4.297 6.685 A 23.013 23.013 0 0 1 3.324 4.316 A 3.831 3.831 0 0 1 3.269 2.047 C 3.567 0.942 4.416 0.2 5.59 0.044 C 6.44 -0.068 7.251 0.015 7.972 0.52 C 9.049 1.274 9.25 2.357 9.125 3.53 C 9.008 4.62 8.6 5.631 8.074 6.595 C 7.623 7.42 7.793 7.944 8.688 8.285 C 9.271 8.506 9.891 8.637 10.495 8.81 C 10.743 8.881 11 8.933 11.235 9.033 C 11.691 9.227 12.028 9.537 12.084 10.051 C 12.144 10.589 12.171 11.436 12.226 12 H 0 Z")'
$ws.Range("G3").Value = 'path("M 7 0 C 10.862 0 14 3.138 14 7 S 10.862 14 7 14 S 0 10.862 0 7 S 3.138 0 7 0 Z M 8.918 9.67 L 8.556 7.518 L 10.111 5.988 L 7.959 5.678 L 7 3.732 L 6.04 5.676 L 3.889 5.987 L 5.444 7.517 L 5.082 9.669 L 7 8.659 L 8.918 9.67 Z")'
$ws.Range("G4").Value = 'path("M 0 8.833 V 0.5 A 0.5 0.5 0 0 1 0.5 0 H 8.833 A 0.5 0.5 0 0 1 9.333 0.5 V 8.833 A 0.5 0.5 0 0 1 8.833 9.333 H 0.5 A 0.5 0.5 0 0 1 0 8.833 Z M 3.468 11.617 V 10.426 H 9.853 A 0.55 0.55 0 0 0 10.425 9.853 V 3.468 H 11.617 V 11.618 H 3.468 Z M 14 14 H 5.851 V 12.809 H 12.236 A 0.55 0.55 0 0 0 12.808 12.236 V 5.851 H 14 V 14 Z")'
$ws.Range("G5").Value = 'path("M 0.215 10.83 C -0.374 9.51 0.288 7.969 1.686 7.404 L 4.42 6.301 L 6.788 11.608 L 6.518 11.718 L 7.508 13.935 L 4.87 15 L 3.873 12.767 C 2.499 13.216 1.004 12.597 0.443 11.339 L 0.215 10.829 Z M 11.132 0.012 C 11.302 0.052 11.442 0.155 11.509 0.306 L 15.95 10.259 C 15.977 10.319 16.004 10.379 16 10.453 C 16.005 10.719 15.765 10.942 15.485 10.947 L 7.135 11.476 L 4.772 6.15 L 10.664 0.148 A 0.507 0.507 0 0 1 11.132 0.013 Z")'
$ws.Range("G6").Value = 'path("M 13.838 0.427 C 13.744 0.229 13.526 0.095 13.272 0.034 C 13.242 0.024 13.228 0.025 13.214 0.023 A 1.194 1.194 0 0 0 12.753 0.023 C 12.705 0.028 12.701 0.032 12.695 0.034 C 12.441 0.095 12.223 0.229 12.13 0.427 L 10.596 9.02 L 8.878 4.047 C 8.784 3.848 8.566 3.714 8.311 3.654 C 8.282 3.644 8.271 3.646 8.258 3.644 A 1.175 1.175 0 0 0 7.791 3.644 C 7.746 3.648 7.742 3.652 7.736 3.654 C 7.482 3.714 7.264 3.848 7.17 4.047 L 5.482 8.502 L 3.726 7.067 C 3.712 7.055 3.691 7.051 3.676 7.04 C 3.624 6.998 3.589 6.992 3.561 6.977 A 1.039 1.039 0 0 0 3.387 6.902 C 3.337 6.887 3.288 6.882 3.236 6.875 C 3.178 6.866 3.122 6.857 3.063 6.857 C 3.003 6.857 2.949 6.867 2.891 6.875 C 2.839 6.883 2.789 6.888 2.74 6.902 A 0.966 0.966 0 0 0 2.566 6.977 C 2.536 6.992 2.503 6.998 2.476 7.017 C 2.434 7.051 2.416 7.055 2.401 7.067 L 0.278 9.382 C -0.048 9.649 -0.101 10.242 0.192 10.555 C 0.486 10.868 0.918 10.821 1.242 10.555 L 3.157 8.545 L 5.236 10.243 C 5.561 10.51 6.121 10.533 6.486 10.293 C 6.524 10.27 6.542 10.238 6.572 10.211 A 0.624 0.624 0 0 0 6.716 10.033 C 6.726 10.017 6.744 10.006 6.752 9.989 L 7.985 6.736 L 10.004 12.574 C 10.097 12.773 10.315 12.905 10.569 12.967 C 10.601 12.977 10.615 12.975 10.629 12.977 A 1.142 1.142 0 0 0 11.087 12.977 C 11.137 12.972 11.14 12.969 11.147 12.967 C 11.4 12.905 11.617 12.772 11.712 12.574 L 13.389 3.176 L 14.429 5.216 C 14.588 5.553 15.234 5.671 15.582 5.473 A 0.862 0.862 0 0 0 15.929 4.396 L 13.838 0.427 Z")'
$ws.Range("G7").Value = 'path("M 11.13 8.524 C 11.47 8.254 11.495 7.896 11.27 7.419 A 19.186 19.186 0 0 1 10.46 5.445 A 3.194 3.194 0 0 1 10.414 3.553 C 10.663 2.631 11.37 2.013 12.348 1.883 C 13.057 1.789 13.733 1.859 14.335 2.28 C 15.233 2.908 15.4 3.811 15.295 4.79 C 15.199 5.698 14.858 6.54 14.419 7.344 C 14.044 8.032 14.185 8.469 14.931 8.753 C 15.418 8.938 15.935 9.047 16.438 9.191 C 16.645 9.251 16.859 9.294 17.055 9.377 C 17.435 9.539 17.716 9.797 17.763 10.225 C 17.813 10.675 17.954 11.53 18 12 H 12.72 C 12.7 11.577 12.681 11.265 12.662 11.064 C 12.623 10.651 12.598 10.264 12.569 9.999 C 12.484 9.22 12.052 8.843 11.229 8.562 A 1.79 1.79 0 0 1 11.13 8.524 Z M 0 12 C 0.09 10.737 0.097 10.387 0.298 9.767 C 0.423 9.381 0.743 9.131 1.135 9.012 C 1.693 8.843 2.265 8.719 2.825 8.557 A 6.36 6.36 0 0 0 3.805 8.207 C 4.508 7.87 4.625 7.382 4.297 6.685 A 23.013 23.013 0 0 1 3.324 4.316 A 3.831 3.831 0 0 1 3.269 2.047 C 3.567 0.942 4.416 0.2 5.59 0.044 C 6.44 -0.068 7.251 0.015 7.972 0.52 C 9.049 1.274 9.25 2.357 9.125 3.53 C 9.008 4.62 8.6 5.631 8.074 6.595 C 7.623 7.42 7.793 7.944 8.688 8.285 C 9.271 8.506 9.891 8.637 10.495 8.81 C 10.743 8.881 11 8.933 11.235 9.033 C 11.691 9.227 12.028 9.537 12.084 10.051 C 12.144 10.589 12.171 11.436 12.226 12 H 0 Z")'

# Apply formatting to match the template row/column styles.
$ws.Range("A2").Copy()
$ws.Range("G3:G6").PasteSpecial(-4122)
$ws.Range("G8:G25").PasteSpecial(-4122)

# Clear the placeholder values in G8:G25 (they only carry formatting, no text)
$ws.Range("G8:G25").ClearContents()

$ws.Range("C8").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Column width for the new column
$ws.Range("G1").ColumnWidth = 32.5703125

# Restore the view: no frozen/top-left offset, selection on A6
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A6").Select()
